$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date in column C for every data row (2..346)
# from 2023-09-23 (45192) to 2023-10-03 (45202).
for ($r = 2; $r -le 346; $r++) {
    $ws.Cells.Item($r, 3).Value = 45202
}

# Row 346 picks up explicit row-height metadata as part of the edit.
$ws.Rows.Item(346).RowHeight = 15

# Append the new row 347 with the new case record.
$ws.Cells.Item(347, 1).Value = "A 46535-2023"

$ws.Cells.Item(347, 2).Value = 45197
$ws.Cells.Item(347, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(347, 3).Value = 45202
$ws.Cells.Item(347, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(347, 4).Value = "VÄSTERBOTTENS LÄN"
$ws.Cells.Item(347, 5).Value = "NORSJÖ"

$ws.Cells.Item(347, 7).Value = 0.8
$ws.Cells.Item(347, 8).Value = 0
$ws.Cells.Item(347, 9).Value = 0
$ws.Cells.Item(347, 10).Value = 0
$ws.Cells.Item(347, 11).Value = 0
$ws.Cells.Item(347, 12).Value = 0
$ws.Cells.Item(347, 13).Value = 0
$ws.Cells.Item(347, 14).Value = 0
$ws.Cells.Item(347, 15).Value = 0
$ws.Cells.Item(347, 16).Value = 0
$ws.Cells.Item(347, 17).Value = 0

$ws.Cells.Item(347, 18).Value = ""
$ws.Cells.Item(347, 18).WrapText = $true
